# Rolling the 2-day forecast window forward by one day:
#  - old rows 98..193 ("day 2", e.g. 24.03.2025) shift up to rows 2..97 ("day 1")
#  - rows 98..193 get a brand new forecast ("day 3", e.g. 25.03.2025)
# Column A = timestamp (date serial), B = forecasted consumption (MW),
# C = quarter index (unchanged), D = "Lookup" label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Snapshot the current "day 2" block (rows 98..193) before overwriting anything ---
$oldRange = $ws.Range("A98:D193")
$old = $oldRange.Value2

# --- 2. New forecast values for the freshly appended day (row 98..193 after the shift) ---
$newB = @(5620,5570,5530,5490,5460,5440,5420,5410,5400,5400,5410,5420,5440,5470,5500,5540,5590,5650,5730,5820,5930,6050,6190,6330,6470,6600,6710,6800,6850,6870,6860,6810,6730,6620,6510,6390,6270,6160,6060,5970,5900,5840,5790,5750,5710,5670,5640,5620,5600,5600,5600,5620,5640,5670,5710,5750,5800,5850,5900,5970,6060,6150,6250,6350,6470,6570,6680,6810,6910,7030,7150,7280,7410,7510,7570,7580,7570,7540,7500,7430,7320,7200,7080,6950,6810,6660,6530,6380,6260,6150,6020,5900,5750,5680,5640,5600)

# --- 3. Build the two write-back blocks ---
$shiftedBlock = $ws.Range("A2:D97").Value2
$newDayBlock  = $ws.Range("A98:D193").Value2

for ($i = 1; $i -le 96; $i++) {
    # Row that used to be "day 2" quarter $i becomes the new "day 1" quarter $i: copy verbatim.
    $shiftedBlock[$i,1] = $old[$i,1]
    $shiftedBlock[$i,2] = $old[$i,2]
    $shiftedBlock[$i,3] = $old[$i,3]
    $shiftedBlock[$i,4] = $old[$i,4]

    # The freshly appended "day 3" quarter $i: same time-of-day one day later, new forecast.
    $newDayBlock[$i,1] = $old[$i,1] + 1
    $newDayBlock[$i,2] = $newB[$i - 1]
    $newDayBlock[$i,3] = $old[$i,3]
    $newDayBlock[$i,4] = "25.03.2025" + $i
}

$ws.Range("A2:D97").Value2   = $shiftedBlock
$ws.Range("A98:D193").Value2 = $newDayBlock
